$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "cost_variable" -> "cost_variable_om" for the data rows C10:C39
$range = $ws.Range("C10:C39")
$range.Value = "cost_variable_om"

# Update the selection to reflect the edited range, matching author's workflow
$range.Select()
